$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new Iran prediction rows (48-51), continuing the pattern of rows 44-47
$ws.Range("A48").Value = "2021-01-09"
$ws.Range("B48").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("D48").Value = 121.54
$ws.Range("F48").Value = "KNN"

$ws.Range("A49").Value = "2021-01-09"
$ws.Range("B49").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D49").Value = 119.89
$ws.Range("F49").Value = "KNN"

$ws.Range("A50").Value = "2021-01-09"
$ws.Range("B50").Value = "24 Jan -- 30 Jan 2021"
$ws.Range("D50").Value = 119.09
$ws.Range("F50").Value = "KNN"

$ws.Range("A51").Value = "2021-01-09"
$ws.Range("B51").Value = "31 Jan -- 06 Feb 2021"
$ws.Range("D51").Value = 115.97
$ws.Range("F51").Value = "KNN"
